# Apply corrected Diebold-Mariano statistics and associated p-values.
$wb = $excel.ActiveWorkbook

# --- Sheet "P_valores": symmetric matrix of p-values ---
$wsP = $wb.Worksheets.Item("P_valores")

$pValues = @{
    "C2" = 0.1517529738683712
    "D2" = 0.4424140926678302
    "E2" = 0.5559930019804442
    "F2" = 0.6335794077072903

    "B3" = 0.1517529738683712
    "D3" = 0.8307987040824378
    "E3" = 0.3568550114688089
    "F3" = 0.5077837894376049

    "B4" = 0.4424140926678302
    "C4" = 0.8307987040824378
    "E4" = 0.5197792387516027
    "F4" = 0.7504395474266516

    "B5" = 0.5559930019804442
    "C5" = 0.3568550114688089
    "D5" = 0.5197792387516027
    "F5" = 0.9486417677235603

    "B6" = 0.6335794077072903
    "C6" = 0.5077837894376049
    "D6" = 0.7504395474266516
    "E6" = 0.9486417677235603
}

foreach ($addr in $pValues.Keys) {
    $wsP.Range($addr).Value = $pValues[$addr]
}

# --- Sheet "Estadisticos_DM": antisymmetric matrix of DM statistics ---
$wsDM = $wb.Worksheets.Item("Estadisticos_DM")

$dmValues = @{
    "C2" = -1.484922795357269
    "D2" = -0.7822365484840662
    "E2" = -0.5979292057262133
    "F2" = -0.4834123394741173

    "B3" = 1.484922795357269
    "D3" = 0.2162327835290688
    "E3" = 0.9411365930245758
    "F3" = 0.6732755844462968

    "B4" = 0.7822365484840662
    "C4" = -0.2162327835290688
    "E4" = 0.6541755710080005
    "F4" = 0.3220707208758357

    "B5" = 0.5979292057262133
    "C5" = -0.9411365930245758
    "D5" = -0.6541755710080005
    "F5" = -0.06515152266377509

    "B6" = 0.4834123394741173
    "C6" = -0.6732755844462968
    "D6" = -0.3220707208758357
    "E6" = 0.06515152266377509
}

foreach ($addr in $dmValues.Keys) {
    $wsDM.Range($addr).Value = $dmValues[$addr]
}
